$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting so values
# like "1.00" or "0.0000270" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "96.693.25"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.720.10"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "239.14"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "1.93"
$ws.Range("E6").Value = "  +9.92%  "
$ws.Range("D7").Value = "655.86"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "0.425"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "1.08"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "3.720.37"
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("D12").Value = "45.41"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "0.206"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  +5.75%  "
$ws.Range("D15").Value = "4.413.06"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "0.0000270"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "96.513.73"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "9.10"
$ws.Range("E18").Value = "  +17.36%  "
$ws.Range("D19").Value = "3.723.71"
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").Value = "19.11"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("D21").Value = "12.88"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "0.532"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "523.34"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "3.50"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "7.06"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "102.89"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "13.40"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  -4.22%  "
$ws.Range("D30").Value = "12.61"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("D31").Value = "3.07"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("E33").Value = "  +12.31%  "
$ws.Range("D34").Value = "0.185"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").Value = "32.86"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "658.49"
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("D38").Value = "0.606"
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("D39").Value = "8.93"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "7.15"
$ws.Range("E40").Value = "  +16.40%  "
$ws.Range("E41").Value = "  +5.00%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "40.36"
$ws.Range("E42").Value = "  +22.68%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.980"
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("D44").Value = "1.99"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "0.457"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "0.0457"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "23.63"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  +1.72%  "
